# Update the "within 100" arithmetic answer table: each table cell holds a
# single run whose text is an expression like "47-9=38". We replace each old
# expression with its new value, one occurrence at a time (wdReplaceOne),
# in document order. A couple of expressions appear twice in the original
# table with two different replacement values, so a single "replace all"
# pass would be wrong for those; doing one match-and-replace per diff
# hunk (in the same top-to-bottom order the hunks appear) reproduces the
# diff exactly, including for the duplicated "50-37=13" cells.
$d = $word.ActiveDocument

$d.Content.Find.Execute("47-9=38", $true, $false, $false, $false, $false, $true, 1, $false, "70-8=62", 1) | Out-Null
$d.Content.Find.Execute("46+26=72", $true, $false, $false, $false, $false, $true, 1, $false, "45-35=10", 1) | Out-Null
$d.Content.Find.Execute("73-55=18", $true, $false, $false, $false, $false, $true, 1, $false, "32+17=49", 1) | Out-Null
$d.Content.Find.Execute("66+23=89", $true, $false, $false, $false, $false, $true, 1, $false, "55+0=55", 1) | Out-Null
$d.Content.Find.Execute("14+34=48", $true, $false, $false, $false, $false, $true, 1, $false, "97-32=65", 1) | Out-Null
$d.Content.Find.Execute("50+40=90", $true, $false, $false, $false, $false, $true, 1, $false, "45+19=64", 1) | Out-Null
$d.Content.Find.Execute("16+22=38", $true, $false, $false, $false, $false, $true, 1, $false, "30+35=65", 1) | Out-Null
$d.Content.Find.Execute("43-6=37", $true, $false, $false, $false, $false, $true, 1, $false, "76-13=63", 1) | Out-Null
$d.Content.Find.Execute("32+39=71", $true, $false, $false, $false, $false, $true, 1, $false, "15+66=81", 1) | Out-Null
$d.Content.Find.Execute("63+26=89", $true, $false, $false, $false, $false, $true, 1, $false, "8+14=22", 1) | Out-Null
$d.Content.Find.Execute("68-31=37", $true, $false, $false, $false, $false, $true, 1, $false, "80-74=6", 1) | Out-Null
$d.Content.Find.Execute("25+72=97", $true, $false, $false, $false, $false, $true, 1, $false, "99-51=48", 1) | Out-Null
$d.Content.Find.Execute("52+24=76", $true, $false, $false, $false, $false, $true, 1, $false, "38+46=84", 1) | Out-Null
$d.Content.Find.Execute("50+28=78", $true, $false, $false, $false, $false, $true, 1, $false, "81-71=10", 1) | Out-Null
$d.Content.Find.Execute("39-31=8", $true, $false, $false, $false, $false, $true, 1, $false, "17+41=58", 1) | Out-Null
$d.Content.Find.Execute("86-7=79", $true, $false, $false, $false, $false, $true, 1, $false, "26+19=45", 1) | Out-Null
$d.Content.Find.Execute("29+62=91", $true, $false, $false, $false, $false, $true, 1, $false, "61+32=93", 1) | Out-Null
$d.Content.Find.Execute("94-37=57", $true, $false, $false, $false, $false, $true, 1, $false, "63+25=88", 1) | Out-Null
$d.Content.Find.Execute("32-20=12", $true, $false, $false, $false, $false, $true, 1, $false, "43-26=17", 1) | Out-Null
$d.Content.Find.Execute("23+72=95", $true, $false, $false, $false, $false, $true, 1, $false, "44+41=85", 1) | Out-Null
$d.Content.Find.Execute("16+13=29", $true, $false, $false, $false, $false, $true, 1, $false, "20+34=54", 1) | Out-Null
$d.Content.Find.Execute("81+4=85", $true, $false, $false, $false, $false, $true, 1, $false, "19+9=28", 1) | Out-Null
$d.Content.Find.Execute("50-42=8", $true, $false, $false, $false, $false, $true, 1, $false, "49+47=96", 1) | Out-Null
$d.Content.Find.Execute("25+24=49", $true, $false, $false, $false, $false, $true, 1, $false, "31+20=51", 1) | Out-Null
$d.Content.Find.Execute("45-21=24", $true, $false, $false, $false, $false, $true, 1, $false, "15+61=76", 1) | Out-Null
$d.Content.Find.Execute("73-63=10", $true, $false, $false, $false, $false, $true, 1, $false, "92-47=45", 1) | Out-Null
$d.Content.Find.Execute("0+49=49", $true, $false, $false, $false, $false, $true, 1, $false, "57+40=97", 1) | Out-Null
$d.Content.Find.Execute("14+66=80", $true, $false, $false, $false, $false, $true, 1, $false, "43-20=23", 1) | Out-Null
$d.Content.Find.Execute("36+31=67", $true, $false, $false, $false, $false, $true, 1, $false, "3+43=46", 1) | Out-Null
$d.Content.Find.Execute("90+6=96", $true, $false, $false, $false, $false, $true, 1, $false, "19+49=68", 1) | Out-Null
$d.Content.Find.Execute("11+25=36", $true, $false, $false, $false, $false, $true, 1, $false, "4+5=9", 1) | Out-Null
$d.Content.Find.Execute("24+43=67", $true, $false, $false, $false, $false, $true, 1, $false, "1+84=85", 1) | Out-Null
$d.Content.Find.Execute("37-27=10", $true, $false, $false, $false, $false, $true, 1, $false, "44-42=2", 1) | Out-Null
$d.Content.Find.Execute("17+48=65", $true, $false, $false, $false, $false, $true, 1, $false, "35-28=7", 1) | Out-Null
$d.Content.Find.Execute("86-24=62", $true, $false, $false, $false, $false, $true, 1, $false, "39+49=88", 1) | Out-Null
$d.Content.Find.Execute("64+24=88", $true, $false, $false, $false, $false, $true, 1, $false, "6-2=4", 1) | Out-Null
$d.Content.Find.Execute("69-20=49", $true, $false, $false, $false, $false, $true, 1, $false, "31+34=65", 1) | Out-Null
$d.Content.Find.Execute("84-75=9", $true, $false, $false, $false, $false, $true, 1, $false, "13+69=82", 1) | Out-Null
$d.Content.Find.Execute("99-96=3", $true, $false, $false, $false, $false, $true, 1, $false, "49-26=23", 1) | Out-Null
$d.Content.Find.Execute("0+87=87", $true, $false, $false, $false, $false, $true, 1, $false, "35-32=3", 1) | Out-Null
$d.Content.Find.Execute("19-16=3", $true, $false, $false, $false, $false, $true, 1, $false, "64-55=9", 1) | Out-Null
$d.Content.Find.Execute("65-64=1", $true, $false, $false, $false, $false, $true, 1, $false, "96-96=0", 1) | Out-Null
$d.Content.Find.Execute("19+11=30", $true, $false, $false, $false, $false, $true, 1, $false, "52-28=24", 1) | Out-Null
$d.Content.Find.Execute("69-18=51", $true, $false, $false, $false, $false, $true, 1, $false, "17+32=49", 1) | Out-Null
$d.Content.Find.Execute("9+13=22", $true, $false, $false, $false, $false, $true, 1, $false, "90-38=52", 1) | Out-Null
$d.Content.Find.Execute("59+20=79", $true, $false, $false, $false, $false, $true, 1, $false, "27+51=78", 1) | Out-Null
$d.Content.Find.Execute("42-0=42", $true, $false, $false, $false, $false, $true, 1, $false, "95-70=25", 1) | Out-Null
$d.Content.Find.Execute("65-25=40", $true, $false, $false, $false, $false, $true, 1, $false, "27+2=29", 1) | Out-Null
$d.Content.Find.Execute("34+42=76", $true, $false, $false, $false, $false, $true, 1, $false, "92+6=98", 1) | Out-Null
$d.Content.Find.Execute("7+68=75", $true, $false, $false, $false, $false, $true, 1, $false, "29+11=40", 1) | Out-Null
$d.Content.Find.Execute("73+0=73", $true, $false, $false, $false, $false, $true, 1, $false, "93+6=99", 1) | Out-Null
$d.Content.Find.Execute("17-14=3", $true, $false, $false, $false, $false, $true, 1, $false, "75+14=89", 1) | Out-Null
$d.Content.Find.Execute("10-6=4", $true, $false, $false, $false, $false, $true, 1, $false, "80-7=73", 1) | Out-Null
$d.Content.Find.Execute("74-44=30", $true, $false, $false, $false, $false, $true, 1, $false, "83-67=16", 1) | Out-Null
$d.Content.Find.Execute("98-4=94", $true, $false, $false, $false, $false, $true, 1, $false, "29-10=19", 1) | Out-Null
$d.Content.Find.Execute("74+13=87", $true, $false, $false, $false, $false, $true, 1, $false, "69+8=77", 1) | Out-Null
$d.Content.Find.Execute("85-10=75", $true, $false, $false, $false, $false, $true, 1, $false, "16+28=44", 1) | Out-Null
$d.Content.Find.Execute("31+65=96", $true, $false, $false, $false, $false, $true, 1, $false, "67-0=67", 1) | Out-Null
$d.Content.Find.Execute("94-24=70", $true, $false, $false, $false, $false, $true, 1, $false, "3+46=49", 1) | Out-Null
$d.Content.Find.Execute("64+14=78", $true, $false, $false, $false, $false, $true, 1, $false, "36+13=49", 1) | Out-Null
$d.Content.Find.Execute("10+52=62", $true, $false, $false, $false, $false, $true, 1, $false, "45+19=64", 1) | Out-Null
$d.Content.Find.Execute("71-30=41", $true, $false, $false, $false, $false, $true, 1, $false, "67-66=1", 1) | Out-Null
$d.Content.Find.Execute("29+38=67", $true, $false, $false, $false, $false, $true, 1, $false, "28+52=80", 1) | Out-Null
$d.Content.Find.Execute("14+41=55", $true, $false, $false, $false, $false, $true, 1, $false, "11+40=51", 1) | Out-Null
$d.Content.Find.Execute("67+3=70", $true, $false, $false, $false, $false, $true, 1, $false, "87-72=15", 1) | Out-Null
$d.Content.Find.Execute("4+90=94", $true, $false, $false, $false, $false, $true, 1, $false, "46-24=22", 1) | Out-Null
$d.Content.Find.Execute("29+55=84", $true, $false, $false, $false, $false, $true, 1, $false, "31-5=26", 1) | Out-Null
$d.Content.Find.Execute("1+78=79", $true, $false, $false, $false, $false, $true, 1, $false, "71-23=48", 1) | Out-Null
$d.Content.Find.Execute("52-38=14", $true, $false, $false, $false, $false, $true, 1, $false, "11-8=3", 1) | Out-Null
$d.Content.Find.Execute("5+45=50", $true, $false, $false, $false, $false, $true, 1, $false, "97-36=61", 1) | Out-Null
$d.Content.Find.Execute("33+33=66", $true, $false, $false, $false, $false, $true, 1, $false, "43-25=18", 1) | Out-Null
$d.Content.Find.Execute("19+40=59", $true, $false, $false, $false, $false, $true, 1, $false, "10+88=98", 1) | Out-Null
$d.Content.Find.Execute("87-63=24", $true, $false, $false, $false, $false, $true, 1, $false, "40+9=49", 1) | Out-Null
$d.Content.Find.Execute("0+46=46", $true, $false, $false, $false, $false, $true, 1, $false, "14+54=68", 1) | Out-Null
$d.Content.Find.Execute("9+53=62", $true, $false, $false, $false, $false, $true, 1, $false, "64-46=18", 1) | Out-Null
$d.Content.Find.Execute("65+11=76", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=11", 1) | Out-Null
$d.Content.Find.Execute("44-39=5", $true, $false, $false, $false, $false, $true, 1, $false, "96-36=60", 1) | Out-Null
$d.Content.Find.Execute("50-37=13", $true, $false, $false, $false, $false, $true, 1, $false, "35-29=6", 1) | Out-Null
$d.Content.Find.Execute("33+63=96", $true, $false, $false, $false, $false, $true, 1, $false, "41-9=32", 1) | Out-Null
$d.Content.Find.Execute("96-31=65", $true, $false, $false, $false, $false, $true, 1, $false, "41-7=34", 1) | Out-Null
$d.Content.Find.Execute("89-39=50", $true, $false, $false, $false, $false, $true, 1, $false, "30-5=25", 1) | Out-Null
$d.Content.Find.Execute("83-3=80", $true, $false, $false, $false, $false, $true, 1, $false, "16+4=20", 1) | Out-Null
$d.Content.Find.Execute("36-7=29", $true, $false, $false, $false, $false, $true, 1, $false, "41+43=84", 1) | Out-Null
$d.Content.Find.Execute("77-73=4", $true, $false, $false, $false, $false, $true, 1, $false, "3+15=18", 1) | Out-Null
$d.Content.Find.Execute("81+2=83", $true, $false, $false, $false, $false, $true, 1, $false, "47+39=86", 1) | Out-Null
$d.Content.Find.Execute("28-4=24", $true, $false, $false, $false, $false, $true, 1, $false, "51-36=15", 1) | Out-Null
$d.Content.Find.Execute("5+47=52", $true, $false, $false, $false, $false, $true, 1, $false, "42+20=62", 1) | Out-Null
$d.Content.Find.Execute("50-37=13", $true, $false, $false, $false, $false, $true, 1, $false, "49+32=81", 1) | Out-Null
$d.Content.Find.Execute("79-64=15", $true, $false, $false, $false, $false, $true, 1, $false, "21-0=21", 1) | Out-Null
$d.Content.Find.Execute("78-44=34", $true, $false, $false, $false, $false, $true, 1, $false, "93-63=30", 1) | Out-Null
$d.Content.Find.Execute("18-15=3", $true, $false, $false, $false, $false, $true, 1, $false, "67-5=62", 1) | Out-Null
$d.Content.Find.Execute("2+71=73", $true, $false, $false, $false, $false, $true, 1, $false, "25-17=8", 1) | Out-Null
$d.Content.Find.Execute("74-21=53", $true, $false, $false, $false, $false, $true, 1, $false, "4+80=84", 1) | Out-Null
$d.Content.Find.Execute("48+36=84", $true, $false, $false, $false, $false, $true, 1, $false, "3+62=65", 1) | Out-Null
$d.Content.Find.Execute("15+11=26", $true, $false, $false, $false, $false, $true, 1, $false, "46-23=23", 1) | Out-Null
$d.Content.Find.Execute("37+7=44", $true, $false, $false, $false, $false, $true, 1, $false, "21+64=85", 1) | Out-Null
$d.Content.Find.Execute("84-27=57", $true, $false, $false, $false, $false, $true, 1, $false, "1+19=20", 1) | Out-Null
$d.Content.Find.Execute("4+17=21", $true, $false, $false, $false, $false, $true, 1, $false, "36+5=41", 1) | Out-Null
$d.Content.Find.Execute("28+41=69", $true, $false, $false, $false, $false, $true, 1, $false, "36-1=35", 1) | Out-Null
$d.Content.Find.Execute("96+3=99", $true, $false, $false, $false, $false, $true, 1, $false, "13+39=52", 1) | Out-Null
